# Refresh the coin price (D) and 1h volume-change (E) columns with the
# latest values from the coinranking.com snapshot used by the GitHub Action.
#
# Price strings that look like plain numbers (e.g. "0.9994") are entered with a
# leading apostrophe so Excel keeps storing them as text, matching the other
# price cells in column D (some of which, like "26.536.21", are not valid
# numbers and are already stored as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '26.536.21'
$ws.Range("E2").Value = '  +4.03%  '
# Row 3 - Ethereum
$ws.Range("D3").Value = '1.739.98'
$ws.Range("E3").Value = '  +4.53%  '
# Row 4 - TetherUSD
$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  -0.02%  '
# Row 5 - BNB
$ws.Range("D5").Value = '''245.56'
$ws.Range("E5").Value = '  +4.80%  '
# Row 6 - USDC
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  -0.07%  '
# Row 7 - XRP
$ws.Range("D7").Value = '''0.4808'
$ws.Range("E7").Value = '  +4.06%  '
# Row 8 - Cardano
$ws.Range("E8").Value = '  +4.52%  '
# Row 9 - Dogecoin
$ws.Range("E9").Value = '  +2.11%  '
# Row 10 - WrappedEther
$ws.Range("D10").Value = '1.739.71'
$ws.Range("E10").Value = '  +4.49%  '
# Row 11 - TRON
$ws.Range("D11").Value = '''0.07128'
$ws.Range("E11").Value = '  +2.79%  '
# Row 12 - Solana
$ws.Range("D12").Value = '''15.85'
$ws.Range("E12").Value = '  +8.65%  '
# Row 13 - Polygon
$ws.Range("D13").Value = '''0.6218'
$ws.Range("E13").Value = '  +8.70%  '
# Row 14 - Polkadot
$ws.Range("D14").Value = '''4.538'
$ws.Range("E14").Value = '  +4.73%  '
# Row 15 - Litecoin
$ws.Range("D15").Value = '''77.12'
$ws.Range("E15").Value = '  +3.07%  '
# Row 16 - Dai
$ws.Range("D16").Value = '''0.9998'
$ws.Range("E16").Value = '  -0.05%  '
# Row 17 - WrappedBTC
$ws.Range("D17").Value = '26.544.76'
$ws.Range("E17").Value = '  +4.01%  '
# Row 18 - BinanceUSD
$ws.Range("E18").Value = '  -0.06%  '
# Row 19 - ShibaInu
$ws.Range("D19").Value = '''0.000006895'
$ws.Range("E19").Value = '  +2.76%  '
# Row 20 - Avalanche
$ws.Range("E20").Value = '  +3.74%  '
# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = '1.961.43'
$ws.Range("E21").Value = '  +4.40%  '
# Row 22 - Uniswap
$ws.Range("D22").Value = '''4.586'
$ws.Range("E22").Value = '  +4.38%  '
# Row 23 - Cosmos
$ws.Range("D23").Value = '''8.910'
$ws.Range("E23").Value = '  +3.16%  '
# Row 24 - Chainlink
$ws.Range("D24").Value = '''5.352'
$ws.Range("E24").Value = '  +2.62%  '
# Row 25 - Monero
$ws.Range("D25").Value = '''135.70'
$ws.Range("E25").Value = '  +0.59%  '
# Row 26 - EthereumClassic
$ws.Range("E26").Value = '  +3.71%  '
# Row 27 - LidoDAOToken
$ws.Range("D27").Value = '''1.813'
$ws.Range("E27").Value = '  +6.12%  '
# Row 28 - Toncoin
$ws.Range("D28").Value = '''1.426'
$ws.Range("E28").Value = '  +4.69%  '
# Row 29 - BitcoinCash
$ws.Range("D29").Value = '''106.88'
$ws.Range("E29").Value = '  +3.31%  '
# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = '''4.011'
$ws.Range("E30").Value = '  +1.63%  '
# Row 31 - Filecoin
$ws.Range("D31").Value = '''3.744'
$ws.Range("E31").Value = '  +4.39%  '
# Row 32 - Stellar
$ws.Range("D32").Value = '''0.07898'
$ws.Range("E32").Value = '  +2.80%  '
# Row 33 - Hedera
$ws.Range("D33").Value = '''0.04589'
$ws.Range("E33").Value = '  +6.58%  '
# Row 34 - HuobiToken
$ws.Range("D34").Value = '''2.614'
$ws.Range("E34").Value = '  -0.25%  '
# Row 35 - ARBITRUM
$ws.Range("D35").Value = '''1.001'
$ws.Range("E35").Value = '  +6.63%  '
# Row 36 - ImmutableX
$ws.Range("D36").Value = '''0.6372'
$ws.Range("E36").Value = '  +6.73%  '
# Row 37 - TrustWalletToken
$ws.Range("D37").Value = '''0.9298'
$ws.Range("E37").Value = '  +1.35%  '
# Row 38 - Quant
$ws.Range("D38").Value = '''111.87'
$ws.Range("E38").Value = '  +5.90%  '
# Row 39 - RenderToken
$ws.Range("D39").Value = '''1.996'
$ws.Range("E39").Value = '  +9.11%  '
# Row 40 - MXToken
$ws.Range("E40").Value = '  -1.76%  '
# Row 41 - PaxDollar
$ws.Range("E41").Value = '  +0.34%  '
# Row 42 - VeChain
$ws.Range("D42").Value = '''0.01516'
$ws.Range("E42").Value = '  +3.99%  '
# Row 43 - FraxShare
$ws.Range("D43").Value = '''5.734'
$ws.Range("E43").Value = '  +15.21%  '
# Row 44 - TheSandbox
$ws.Range("D44").Value = '''0.3920'
$ws.Range("E44").Value = '  +6.00%  '
# Row 45 - Aptos
$ws.Range("D45").Value = '''6.961'
$ws.Range("E45").Value = '  +14.31%  '
# Row 46 - Algorand
$ws.Range("D46").Value = '''0.1199'
$ws.Range("E46").Value = '  +8.08%  '
# Row 47 - Cronos
$ws.Range("D47").Value = '''0.05334'
$ws.Range("E47").Value = '  +1.54%  '
# Row 48 - EnergySwap
$ws.Range("D48").Value = '''7.905'
$ws.Range("E48").Value = '  +5.22%  '
# Row 49 - Elrond
$ws.Range("D49").Value = '''30.89'
$ws.Range("E49").Value = '  +3.36%  '
# Row 50 - NEARProtocol
$ws.Range("D50").Value = '''1.259'
$ws.Range("E50").Value = '  +6.45%  '
# Row 51 - Decentraland
$ws.Range("D51").Value = '''0.3454'
$ws.Range("E51").Value = '  +4.92%  '
